$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("meta")

# Row 5 previously held just an (empty) styled placeholder cell in column A.
# Turn it into a real "style" / "default" key-value pair, and push the
# empty styled placeholder down to the new row 6.
$a4 = $meta.Range("A4")
$a5 = $meta.Range("A5")
$a6 = $meta.Range("A6")

# Copy the bold/orange label formatting from A4 onto the new A5 and A6 cells.
$a4.Copy()
$a5.PasteSpecial(-4122)   # xlPasteFormats
$a6.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$a5.Value = "style"
$meta.Range("B5").Value = "default"
